# Iniciado desenvolvimento da barra de pesquisa no Frontend
# Reworks the "ok" markers in column F of Planilha1:
#  - F4's marker is removed (cleared)
#  - F5 gains a new "ok" marker (new cell/format)
#  - F9 and F19 flip their marker text from "ok" to "o"
#  - F14, F15 and F21 gain new "o" markers
#  - F17 and F18 have their "ok" markers cleared (cell/style stays)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F4: remove the existing "ok" marker entirely.
$ws.Range("F4").ClearContents()

# F5: new "ok" marker (vertically centered, matching the column's style).
$ws.Range("F5").Value = "ok"
$ws.Range("F5").VerticalAlignment = -4108

# F9: marker text changes from "ok" to "o".
$ws.Range("F9").Value = "o"

# F14: new "o" marker.
$ws.Range("F14").Value = "o"
$ws.Range("F14").VerticalAlignment = -4108

# F15: new "o" marker.
$ws.Range("F15").Value = "o"
$ws.Range("F15").VerticalAlignment = -4108

# F17 and F18: clear the "ok" markers, keeping their existing style.
$ws.Range("F17").ClearContents()
$ws.Range("F18").ClearContents()

# F19: marker text changes from "ok" to "o".
$ws.Range("F19").Value = "o"

# F21: new "o" marker.
$ws.Range("F21").Value = "o"
$ws.Range("F21").VerticalAlignment = -4108

# Update the sheet's remembered selection.
$ws.Range("K17").Select()
